$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, shifting existing rows 27-112 down to 28-113.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new record.
$ws.Range("A27").Value = 5
$ws.Range("B27").Value = "Macroferia Regional de Talca"
$ws.Range("C27").Value = "Maule"
$ws.Range("D27").Value = 45076
$ws.Range("E27").Value = 7
$ws.Range("F27").Value = 100112013
$ws.Range("G27").Value = "Alcachofa"
$ws.Range("H27").Value = "Española"
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 300
$ws.Range("K27").Value = 18000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 18000
$ws.Range("N27").Value = "`$/caja 40 unidades"
$ws.Range("O27").Value = "Provincia de Limarí"
$ws.Range("P27").Value = 450
$ws.Range("Q27").Value = 40
$ws.Range("R27").Value = "Hortaliza"
